$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty monthly data columns (G:J and L:N) for rows 3-16.
$data = @{
    3  = @(1502, 1371, 1544, 1676, 1433, 1534, 1679)
    4  = @(6497, 5868, 6961, 8361, 6152, 5853, 6398)
    5  = @(480,  408,  476,  512,  479,  487,  585)
    6  = @(120,  97,   116,  122,  94,   88,   72)
    7  = @(12219,11567,13005,13978,12823,12975,13799)
    8  = @(2913, 2678, 3387, 3702, 3935, 4137, 4194)
    9  = @(8431, 7720, 7994, 8737, 6584, 6910, 9023)
    10 = @(2340, 2191, 2544, 3563, 2445, 2464, 2827)
    11 = @(400,  386,  426,  542,  382,  404,  406)
    12 = @(19585,19969,20731,22158,18293,17388,18775)
    13 = @(1896, 1908, 2359, 2663, 2391, 2465, 2948)
    14 = @(27276,24852,2811, 32663,24537,24818,26996)
    15 = @(2336, 2053, 2243, 2455, 2309, 2419, 2517)
    16 = @(9868, 9875, 10546,10747,9532, 9619, 10693)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("G$row").Value = $vals[0]
    $ws.Range("H$row").Value = $vals[1]
    $ws.Range("I$row").Value = $vals[2]
    $ws.Range("J$row").Value = $vals[3]
    $ws.Range("L$row").Value = $vals[4]
    $ws.Range("M$row").Value = $vals[5]
    $ws.Range("N$row").Value = $vals[6]
}

# Update the selected range/active cell to match the saved view state.
$ws.Range("L17").Select() | Out-Null
